$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:tab/></w:r><w:r><w:t>Everyday</w:t></w:r><w:r><w:t xml:space="preserve"> millions of users have to decide what</w:t></w:r><w:r><w:t xml:space="preserve"> they want to eat for breakfast, lunch, and dinner. For various reasons</w:t></w:r><w:r><w:t>—</w:t></w:r><w:r><w:t>including the amount of dining options, dietary concerns, and price range—many of those people have a hard time deciding where to spend their money.</w:t></w:r><w:r><w:t xml:space="preserve"> While apps like </w:t></w:r><w:r><w:t>G</w:t></w:r><w:r><w:t xml:space="preserve">oogle </w:t></w:r><w:r><w:t>M</w:t></w:r><w:r><w:t>aps are likely to overwhelm users with a long list of options, our app will help people filter through the options.</w:t></w:r><w:r><w:t xml:space="preserve"> To help </w:t></w:r><w:r><w:t>these</w:t></w:r><w:r><w:t xml:space="preserve"> users decide on dinning options in a timely manner we</w:t></w:r><w:r><w:t xml:space="preserve"> will create an app that can help make decisions based on preferences, distance, and price-range.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r></w:p><w:p><w:r><w:tab/><w:t>Anyone could use this app to help them, but</w:t></w:r><w:r><w:t xml:space="preserve"> the app would mostly be targeted at people who frequently eat out, especially those who eat a lot of fast food. Because of this the age demographic will most likely be young around student age.</w:t></w:r></w:p><w:p/><w:p><w:r><w:tab/><w:t xml:space="preserve"> At the beginning a user sign-in page will be presented. Each unique user will sign in to see their own preferences and suggestions, and those will be saved locally for when the user returns. Next the user will be presented several buttons. One button will be to see all the restaurants sorted by distance in a certain range. Another will be to see only favorited or frequently visited locations. Next </w:t></w:r><w:r><w:t>a button that will quiz a user on what kind of foods they are in the mood for, what their price range is, etc. Lastly the user can press a “random” button that will randomly suggest a restaurant and show them (they will be able to hit this button until they find something they want). The app will also display the relative busyness, hours of operation, and estimate price of the restaurant.</w:t></w:r></w:p><w:p/><w:p><w:r><w:tab/><w:t xml:space="preserve">For this project we will have to borrow heavily from the Google Maps API. We will also use some method of </w:t></w:r><w:r><w:t>encryption</w:t></w:r><w:r><w:t xml:space="preserve"> to hold the usernames and passwords. We would also like to look at databases to store user preferences.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$target = $d.Content
$target.Find.Execute("Write some stuff")
$target.InsertXML($xml)
